$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.926.61"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "3.950.35"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "3.941.08"
$ws.Range("E7").Value = "  -2.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.689"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -6.34%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.737"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.34"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +15.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000317"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.13%  "
$ws.Range("D15").Value = "4.581.16"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "3.952.33"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.74"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -4.11%  "
$ws.Range("D21").Value = "70.898.45"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.87"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "96.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.15"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +15.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.37"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +12.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.70"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +14.60%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "690.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.29"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.129"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "64.19"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.435"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.48"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("D40").Value = "0.0₃0818"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.23"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0480"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.148"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.71"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.40"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.98"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000275"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.31%  "
